$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-15 22:48:53'
$ws.Range('E3').Value = '2026-02-15 22:48:56'
$ws.Range('I3').Value = '3.1 mm'
$ws.Range('O3').Value = '-4.7 °C'
$ws.Range('E4').Value = '2026-02-15 22:48:59'
$ws.Range('E5').Value = '2026-02-15 22:49:01'
$ws.Range('I5').Value = '9.2 mm'
$ws.Range('O5').Value = '-4.2 °C'
$ws.Range('E6').Value = '2026-02-15 22:49:04'
$ws.Range('E7').Value = '2026-02-15 22:49:07'
$ws.Range('O7').Value = '12.0 °C'
$ws.Range('E8').Value = '2026-02-15 22:49:10'
$ws.Range('O8').Value = '8.4 °C'
$ws.Range('E9').Value = '2026-02-15 22:49:12'
$ws.Range('O9').Value = '10.6 °C'
$ws.Range('E10').Value = '2026-02-15 22:49:15'
$ws.Range('E11').Value = '2026-02-15 22:49:18'
$ws.Range('H11').Value = "'48%"
$ws.Range('O11').Value = '6.8 °C'
$ws.Range('E12').Value = '2026-02-15 22:49:20'
$ws.Range('N12').Value = '7.9 °C 22:16 TU'
$ws.Range('O12').Value = '10.6 °C'
$ws.Range('E13').Value = '2026-02-15 22:49:23'
$ws.Range('H13').Value = "'41%"
$ws.Range('J13').Value = '1015.6 hPa'
$ws.Range('O13').Value = '6.2 °C'
$ws.Range('E14').Value = '2026-02-15 22:49:26'
$ws.Range('E15').Value = '2026-02-15 22:49:28'
$ws.Range('H15').Value = "'54%"
$ws.Range('O15').Value = '10.3 °C'
$ws.Range('E16').Value = '2026-02-15 22:49:31'
$ws.Range('H16').Value = "'64%"
$ws.Range('I16').Value = '1.3 mm'
$ws.Range('E17').Value = '2026-02-15 22:49:34'
$ws.Range('O17').Value = '3.2 °C'
$ws.Range('E18').Value = '2026-02-15 22:49:36'
$ws.Range('E19').Value = '2026-02-15 22:49:39'
$ws.Range('E20').Value = '2026-02-15 22:49:41'
$ws.Range('E21').Value = '2026-02-15 22:49:44'
$ws.Range('E22').Value = '2026-02-15 22:49:47'
$ws.Range('N22').Value = '-6.5 °C 22:15 TU'
$ws.Range('E23').Value = '2026-02-15 22:49:49'
$ws.Range('H23').Value = "'68%"
$ws.Range('I23').Value = '6.0 mm'
$ws.Range('O23').Value = '-3.3 °C'
$ws.Range('E24').Value = '2026-02-15 22:49:52'
$ws.Range('O24').Value = '9.1 °C'
$ws.Range('E25').Value = '2026-02-15 22:49:55'
$ws.Range('H25').Value = "'64%"
$ws.Range('O25').Value = '-1.2 °C'
$ws.Range('E26').Value = '2026-02-15 22:49:57'
$ws.Range('E27').Value = '2026-02-15 22:50:00'
$ws.Range('H27').Value = "'52%"
$ws.Range('E28').Value = '2026-02-15 22:50:02'
$ws.Range('E29').Value = '2026-02-15 22:50:05'
$ws.Range('H29').Value = "'61%"
$ws.Range('O29').Value = '9.9 °C'
$ws.Range('E30').Value = '2026-02-15 22:50:07'
$ws.Range('E31').Value = '2026-02-15 22:50:10'
$ws.Range('E32').Value = '2026-02-15 22:50:13'
$ws.Range('E33').Value = '2026-02-15 22:50:15'
$ws.Range('O33').Value = '6.1 °C'
$ws.Range('E34').Value = '2026-02-15 22:50:18'
$ws.Range('E35').Value = '2026-02-15 22:50:20'
$ws.Range('J35').Value = '1019.6 hPa'
$ws.Range('E36').Value = '2026-02-15 22:50:23'
$ws.Range('N36').Value = '7.6 °C 22:19 TU'
$ws.Range('O36').Value = '11.2 °C'
$ws.Range('E37').Value = '2026-02-15 22:50:26'
$ws.Range('H37').Value = "'57%"
$ws.Range('J37').Value = '1016.6 hPa'
$ws.Range('O37').Value = '5.7 °C'
$ws.Range('E38').Value = '2026-02-15 22:50:28'
$ws.Range('E39').Value = '2026-02-15 22:50:31'
$ws.Range('H39').Value = "'61%"
$ws.Range('O39').Value = '-2.5 °C'
$ws.Range('E40').Value = '2026-02-15 22:50:33'
$ws.Range('H40').Value = "'43%"
$ws.Range('E41').Value = '2026-02-15 22:50:36'
$ws.Range('E42').Value = '2026-02-15 22:50:38'
$ws.Range('H42').Value = "'60%"
$ws.Range('O42').Value = '10.4 °C'
$ws.Range('E43').Value = '2026-02-15 22:50:41'
$ws.Range('E44').Value = '2026-02-15 22:50:44'
$ws.Range('H44').Value = "'79%"
$ws.Range('I44').Value = '6.5 mm'
$ws.Range('E45').Value = '2026-02-15 22:50:47'
$ws.Range('H45').Value = "'91%"
$ws.Range('I45').Value = '4.9 mm'
$ws.Range('E46').Value = '2026-02-15 22:50:49'
$ws.Range('H46').Value = "'53%"
